$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 / column B holds the rule name for the "R40" rule ("Good Night").
# Rename it from "R40" to "1" (kept as text, matching the other rule-name
# cells in this column which are all stored as strings).
$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
